# Apply "configurable exam duration" change:
#   - All exam durations are normalized to a 2-hour (120 minute) session.
#   - Morning slot "09:00 - 12:00" becomes "09:00 - 11:00".
#   - Afternoon slot "14:00 - 17:00" becomes "14:00 - 16:00".
#   - Dependent configuration / summary sheets are recomputed to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Exam_Schedule  (columns: D=duration, E=duration_minutes, K=time_slot)
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Exam_Schedule")
$lastRow = $schedule.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $schedule.Cells.Item($r, 4).Value2 = "2 hours"   # duration
    $schedule.Cells.Item($r, 5).Value2 = 120          # duration_minutes

    $slotCell = $schedule.Cells.Item($r, 11)          # time_slot
    $slotVal = $slotCell.Value2
    if ($slotVal -eq "09:00 - 12:00") {
        $slotCell.Value2 = "09:00 - 11:00"
    } elseif ($slotVal -eq "14:00 - 17:00") {
        $slotCell.Value2 = "14:00 - 16:00"
    }
}

# ---------------------------------------------------------------------------
# Sheet: Exam_Classrooms  (columns: E=Time Slot, I=Duration)
# ---------------------------------------------------------------------------
$classrooms = $wb.Worksheets.Item("Exam_Classrooms")
$lastRow2 = $classrooms.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow2; $r++) {
    $slotCell = $classrooms.Cells.Item($r, 5)         # Time Slot
    $slotVal = $slotCell.Value2
    if ($slotVal -eq "09:00 - 12:00") {
        $slotCell.Value2 = "09:00 - 11:00"
    } elseif ($slotVal -eq "14:00 - 17:00") {
        $slotCell.Value2 = "14:00 - 16:00"
    }

    $classrooms.Cells.Item($r, 9).Value2 = "2 hours"  # Duration
}

# ---------------------------------------------------------------------------
# Sheet: Configuration (Session Duration (minutes) parameter)
# ---------------------------------------------------------------------------
$config = $wb.Worksheets.Item("Configuration")
$configLastRow = $config.UsedRange.Rows.Count
for ($r = 2; $r -le $configLastRow; $r++) {
    $label = $config.Cells.Item($r, 1).Value2
    if ($label -eq "Session Duration (minutes)") {
        $config.Cells.Item($r, 2).Value2 = 120
    }
}

# ---------------------------------------------------------------------------
# Sheet: Department_Summary (Total Duration columns recomputed with the new
# uniform 120 minute session length for every exam)
# ---------------------------------------------------------------------------
$deptSummary = $wb.Worksheets.Item("Department_Summary")
$deptLastRow = $deptSummary.UsedRange.Rows.Count

for ($r = 2; $r -le $deptLastRow; $r++) {
    $examCount = $deptSummary.Cells.Item($r, 2).Value2
    $totalMinutes = [double]$examCount * 120
    $deptSummary.Cells.Item($r, 3).Value2 = $totalMinutes
    $deptSummary.Cells.Item($r, 5).Value2 = $totalMinutes / 60
}

Write-Host "Applied configurable exam duration update."
